$wb = $excel.ActiveWorkbook

# Map of sheet index -> worksheet object
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 ---
$ws1.Cells.Item(3, 6).Value = 1211
$ws1.Cells.Item(4, 6).Value = 12738
$ws1.Cells.Item(5, 6).Value = 720
$ws1.Cells.Item(7, 6).Value = 317
$ws1.Cells.Item(8, 6).Value = 55
$ws1.Cells.Item(10, 6).Value = 1861
$ws1.Cells.Item(13, 6).Value = 512
$ws1.Cells.Item(14, 6).Value = 207
$ws1.Cells.Item(16, 6).Value = 339
$ws1.Cells.Item(18, 6).Value = 293
$ws1.Cells.Item(19, 6).Value = 128
$ws1.Cells.Item(20, 6).Value = 124
$ws1.Cells.Item(23, 6).Value = 241
$ws1.Cells.Item(24, 6).Value = 1277
$ws1.Cells.Item(25, 6).Value = 330
$ws1.Cells.Item(27, 6).Value = 102
# --- 演出 ---
$ws2.Cells.Item(4, 6).Value = 281
$ws2.Cells.Item(6, 6).Value = 145
$ws2.Cells.Item(11, 6).Value = 352
$ws2.Cells.Item(17, 6).Value = 9
$ws2.Cells.Item(19, 6).Value = 13
# --- 本地生活 ---
$ws3.Cells.Item(2, 6).Value = 869
# --- 全部类型 ---
$ws4.Cells.Item(2, 6).Value = 869
$ws4.Cells.Item(6, 6).Value = 1211
$ws4.Cells.Item(7, 6).Value = 12738
$ws4.Cells.Item(8, 6).Value = 281
$ws4.Cells.Item(9, 6).Value = 720
$ws4.Cells.Item(12, 6).Value = 317
$ws4.Cells.Item(13, 6).Value = 55
$ws4.Cells.Item(15, 6).Value = 1861
$ws4.Cells.Item(18, 6).Value = 512
$ws4.Cells.Item(20, 6).Value = 207
$ws4.Cells.Item(21, 6).Value = 145
$ws4.Cells.Item(22, 6).Value = 145
$ws4.Cells.Item(28, 6).Value = 352
$ws4.Cells.Item(29, 6).Value = 339
$ws4.Cells.Item(32, 6).Value = 293
$ws4.Cells.Item(33, 6).Value = 128
$ws4.Cells.Item(34, 6).Value = 124
$ws4.Cells.Item(40, 6).Value = 241
$ws4.Cells.Item(41, 6).Value = 1277
$ws4.Cells.Item(43, 6).Value = 330
$ws4.Cells.Item(45, 6).Value = 102
$ws4.Cells.Item(46, 6).Value = 9
$ws4.Cells.Item(48, 6).Value = 13

# Special updates: F and G columns change together (item temporarily sold out)
$ws3.Cells.Item(3, 6).Value = 3642
$ws3.Cells.Item(3, 7).Value = "暂时售罄"
$ws4.Cells.Item(10, 6).Value = 3642
$ws4.Cells.Item(10, 7).Value = "暂时售罄"
